# This script applies the updated NATMI TPM-derived ligand-receptor
# statistics (commit: "update scripts wuth new tpm") to Sheet1, rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.3454506666666666
$arr[0,3] = 1.036352
$arr[0,4] = 0.1052716477644991
$arr[0,5] = 0.1052716477644991
$ws.Range("E2:J2").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 1.845238666666667
$arr[0,1] = 5.535716000000001
$arr[0,2] = 0.1139123054545619
$arr[0,3] = 0.1139123054545619
$arr[0,4] = 0.6374389275591111
$arr[0,5] = 5.736950348032001
$arr[0,6] = 0.01199173609585468
$arr[0,7] = 0.01199173609585468
$ws.Range("M2:T2").Value = $arr

# Row 3
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.3454506666666666
$arr[0,3] = 1.036352
$arr[0,4] = 0.1052716477644991
$arr[0,5] = 0.1052716477644991
$ws.Range("E3:J3").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.6403884287886557
$arr[0,1] = 0.6403884287886557
$arr[0,2] = 3.583533066417778
$arr[0,3] = 32.25179759776
$arr[0,4] = 0.0674147451079004
$arr[0,5] = 0.0674147451079004
$ws.Range("O3:T3").Value = $arr

# Row 4
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3
$arr[0,1] = 1
$arr[0,2] = 0.3454506666666666
$arr[0,3] = 1.036352
$arr[0,4] = 0.1052716477644991
$arr[0,5] = 0.1052716477644991
$ws.Range("E4:J4").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.2456992657567824
$arr[0,1] = 0.2456992657567825
$arr[0,2] = 1.374902174449778
$arr[0,3] = 12.374119570048
$arr[0,4] = 0.02586516656074406
$arr[0,5] = 0.02586516656074407
$ws.Range("O4:T4").Value = $arr

# Row 5
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3398937483175971
$arr[0,1] = 0.3398937483175971
$ws.Range("I5:J5").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 1.845238666666667
$arr[0,1] = 5.535716000000001
$arr[0,2] = 0.1139123054545619
$arr[0,3] = 0.1139123054545619
$arr[0,4] = 2.058118315924
$arr[0,5] = 18.523064843316
$arr[0,6] = 0.03871808048045011
$arr[0,7] = 0.03871808048045011
$ws.Range("M5:T5").Value = $arr

# Row 6
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3398937483175971
$arr[0,1] = 0.3398937483175971
$ws.Range("I6:J6").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.6403884287886557
$arr[0,1] = 0.6403884287886557
$ws.Range("O6:P6").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2176640234401928
$arr[0,1] = 0.2176640234401928
$ws.Range("S6:T6").Value = $arr

# Row 7
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3398937483175971
$arr[0,1] = 0.3398937483175971
$ws.Range("I7:J7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2456992657567824
$arr[0,1] = 0.2456992657567825
$ws.Range("O7:P7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.0835116443969542
$arr[0,1] = 0.08351164439695422
$ws.Range("S7:T7").Value = $arr

# Row 8
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5548346039179038
$arr[0,1] = 0.5548346039179038
$ws.Range("I8:J8").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 1.845238666666667
$arr[0,1] = 5.535716000000001
$arr[0,2] = 0.1139123054545619
$arr[0,3] = 0.1139123054545619
$arr[0,4] = 3.359624195161334
$arr[0,5] = 30.236617756452
$arr[0,6] = 0.06320248887825713
$arr[0,7] = 0.06320248887825715
$ws.Range("M8:T8").Value = $arr

# Row 9
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5548346039179038
$arr[0,1] = 0.5548346039179038
$ws.Range("I9:J9").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.6403884287886557
$arr[0,1] = 0.6403884287886557
$ws.Range("O9:P9").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3553096602405625
$arr[0,1] = 0.3553096602405625
$ws.Range("S9:T9").Value = $arr

# Row 10
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.5548346039179038
$arr[0,1] = 0.5548346039179038
$ws.Range("I10:J10").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.2456992657567824
$arr[0,1] = 0.2456992657567825
$ws.Range("O10:P10").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.1363224547990842
$arr[0,1] = 0.1363224547990842
$ws.Range("S10:T10").Value = $arr
